# Secondary screening results update
# Updates LLM-generated summary/analysis text and associated extracted
# fields for several rows in the secondary-screening results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 'I am sorry, but I cannot complete your request as you have not provided the article text for analysis.'

# --- Row 3 (previously a LangFlow error row; now populated with a real, if empty-input, response) ---
$ws.Range("E3").Value = 'I am sorry, but I am missing the article text. Please provide it so that I can complete the analysis.'
$ws.Range("F3").Value = 'N/A'
$ws.Range("G3").Value = 'N/A'
$ws.Range("H3").Value = 'N/A'
$ws.Range("I3").Value = 'N/A'
$ws.Range("J3").Value = 'N/A'
$ws.Range("K3").Value = 'N/A'
$ws.Range("L3").Value = 'N/A'
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 'T2'
$ws.Range("O3").Value = 'O1'
$ws.Range("P3").Value = 'F2'
$ws.Range("Q3").Value = 'S2'
$ws.Range("R3").Value = 'C1'
$ws.Range("S3").Value = 8
$ws.Range("T3").Value = 'NA'
$ws.Range("U3").Value = 'NA'
$ws.Range("V3").Value = 'NA'
$ws.Range("W3").Value = 'include'

# --- Row 4 ---
$ws.Range("E4").Value = 'The provided text is empty, so no summary can be generated.'
$ws.Range("F4").Value = 'Cannot be determined from empty text.'
$ws.Range("G4").Value = 'Cannot be determined from empty text.'
$ws.Range("H4").Value = 'Cannot be determined from empty text.'
$ws.Range("I4").Value = 'Cannot be determined from empty text.'
$ws.Range("J4").Value = 'Cannot be determined from empty text.'
$ws.Range("K4").Value = 'Cannot be determined from empty text.'
$ws.Range("L4").Value = 'Cannot be determined from empty text.'

# --- Row 5 ---
$ws.Range("E5").Value = 'I am sorry, but I am missing the article text. I need the article text to complete the JSON.'

# --- Row 6 (previously a LangFlow error row; now populated with a real, if empty-input, response) ---
$ws.Range("E6").Value = 'I am sorry, but I am missing the article text. Please provide the article text so that I can properly analyze it and provide the appropriate JSON output.'
$ws.Range("F6").Value = 'N/A'
$ws.Range("G6").Value = 'N/A'
$ws.Range("H6").Value = 'N/A'
$ws.Range("I6").Value = 'N/A'
$ws.Range("J6").Value = 'N/A'
$ws.Range("K6").Value = 'N/A'
$ws.Range("L6").Value = 'N/A'
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 'T2'
$ws.Range("O6").Value = 'O1'
$ws.Range("P6").Value = 'F2'
$ws.Range("Q6").Value = 'S2'
$ws.Range("R6").Value = 'C1'
$ws.Range("S6").Value = 8
$ws.Range("T6").Value = 'NA'
$ws.Range("U6").Value = 'NA'
$ws.Range("V6").Value = 'NA'
$ws.Range("W6").Value = 'include'

# --- Row 7 ---
$ws.Range("E7").Value = 'I am sorry, but I do not have access to either the IFU context or the article text. Therefore, I am unable to provide a JSON response based on the article.'

# --- Row 8 ---
$ws.Range("E8").Value = 'I am sorry, but I am missing the article text to provide an accurate response.'

# --- Row 9 ---
$ws.Range("E9").Value = 'I am sorry, but I need the content of the article and IFU to provide a complete response.'

# --- Row 10 ---
$ws.Range("E10").Value = 'I am sorry, but I am missing the article and IFU context. Please provide the data so that I can answer the question.'

# --- Row 11 (now becomes a LangFlow error row, clearing the previously extracted fields) ---
$ws.Range("E11").Value = 'LangFlow error: Expecting value: line 1 column 1 (char 0)'
$ws.Range("F11:V11").ClearContents()
$ws.Range("W11").Value = 'exclude'

# --- Row 12 ---
$ws.Range("E12").Value = 'I am sorry, but I am missing the article text. Please provide the article text so that I can provide an accurate response.'
$ws.Range("F12").Value = 'N/A'
$ws.Range("G12").Value = 'N/A'
$ws.Range("H12").Value = 'N/A'
$ws.Range("I12").Value = 'N/A'
$ws.Range("J12").Value = 'N/A'
$ws.Range("K12").Value = 'N/A'
$ws.Range("L12").Value = 'N/A'

# --- Row 13 ---
$ws.Range("E13").Value = 'I am sorry, but I need the article text to provide a JSON formatted response.'

# --- Row 14 ---
$ws.Range("E14").Value = 'I am sorry, but I cannot process your request because the article text is missing.'

Write-Host "Secondary screening results updated."
